$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 181, shifting existing rows (181-265) down to (182-266)
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new record
$ws.Cells.Item(181,1).Value2 = 5
$ws.Cells.Item(181,2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(181,3).Value2 = "Maule"
$ws.Cells.Item(181,4).Value2 = 44917
$ws.Cells.Item(181,4).NumberFormat = $ws.Cells.Item(182,4).NumberFormat
$ws.Cells.Item(181,5).Value2 = 7
$ws.Cells.Item(181,6).Value2 = 100112024
$ws.Cells.Item(181,7).Value2 = "Choclo"
$ws.Cells.Item(181,8).Value2 = "Choclero"
$ws.Cells.Item(181,9).Value2 = "Primera"
$ws.Cells.Item(181,10).Value2 = 20000
$ws.Cells.Item(181,11).Value2 = 300
$ws.Cells.Item(181,12).Value2 = 300
$ws.Cells.Item(181,13).Value2 = 300
$ws.Cells.Item(181,14).Value2 = "$/unidad"
$ws.Cells.Item(181,15).Value2 = "Región del Maule"
$ws.Cells.Item(181,16).Value2 = 300
$ws.Cells.Item(181,17).Value2 = 1
$ws.Cells.Item(181,18).Value2 = "Hortaliza"
